$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data row (2..360) has its "Förändrad" (column C) date bumped
# from 45204 (2023-10-05) to 45205 (2023-10-06).
$ws.Range("C2:C360").Value = 45205

# Row 3 (A 19109-2022) additionally gained a new signal species
# (Motaggsvamp) and a new red-listed species (Grönpyrola), which bumps
# several summary counts and updates the species list text.
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 9
$ws.Range("O3").Value = 11
$ws.Range("Q3").Value = 18
$ws.Range("R3").Value = "Knärot`r`nRynkskinn`r`nGammelgransskål`r`nGarnlav`r`nGranticka`r`nJärpe`r`nLunglav`r`nMotaggsvamp`r`nTretåig hackspett`r`nUllticka`r`nViolettgrå tagellav`r`nBollvitmossa`r`nBårdlav`r`nGrönpyrola`r`nPlattlummer`r`nSpindelblomster`r`nStuplav`r`nVedticka"

# Row 5 (A 21627-2023) additionally gained a new signal species
# (Dropptaggsvamp) and a new species (Vedticka), which bumps the
# signal-species and all-species counts and updates the species list text.
$ws.Range("I5").Value = 7
$ws.Range("Q5").Value = 14
$ws.Range("R5").Value = "Knärot`r`nRynkskinn`r`nGarnlav`r`nGranticka`r`nLunglav`r`nUllticka`r`nVitgrynig nållav`r`nBårdlav`r`nDropptaggsvamp`r`nLuddlav`r`nSkinnlav`r`nStor aspticka`r`nStuplav`r`nVedticka"
